$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41; this shifts existing rows 41-95 down to 42-96
$ws.Rows.Item(41).Insert()

# Copy static/unchanged column values from the row that is now directly below (row 42,
# which holds what used to be row 41's data) into the newly inserted row 41.
$ws.Cells.Item(41, 1).Value2 = $ws.Cells.Item(42, 1).Value2   # A Mercado ID
$ws.Cells.Item(41, 2).Value2 = $ws.Cells.Item(42, 2).Value2   # B Mercado
$ws.Cells.Item(41, 3).Value2 = $ws.Cells.Item(42, 3).Value2   # C Region
$ws.Cells.Item(41, 5).Value2 = $ws.Cells.Item(42, 5).Value2   # E Codreg
$ws.Cells.Item(41, 6).Value2 = $ws.Cells.Item(42, 6).Value2   # F Categoria ID
$ws.Cells.Item(41, 7).Value2 = $ws.Cells.Item(42, 7).Value2   # G Categoria
$ws.Cells.Item(41, 8).Value2 = $ws.Cells.Item(42, 8).Value2   # H Variedad
$ws.Cells.Item(41, 9).Value2 = $ws.Cells.Item(42, 9).Value2   # I Calidad
$ws.Cells.Item(41, 14).Value2 = $ws.Cells.Item(42, 14).Value2 # N Unidad de comercializacion
$ws.Cells.Item(41, 15).Value2 = $ws.Cells.Item(42, 15).Value2 # O Origen
$ws.Cells.Item(41, 17).Value2 = $ws.Cells.Item(42, 17).Value2 # Q Kg o Unidades
$ws.Cells.Item(41, 18).Value2 = $ws.Cells.Item(42, 18).Value2 # R Clasificacion

# Match the number format of the date column (D) for the new row, same as the row below it.
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat

# New row's own data values
$ws.Cells.Item(41, 4).Value2 = 44895   # D Fecha
$ws.Cells.Item(41, 10).Value2 = 250    # J Volumen
$ws.Cells.Item(41, 11).Value2 = 1500   # K Precio minimo
$ws.Cells.Item(41, 12).Value2 = 1800   # L Precio maximo
$ws.Cells.Item(41, 13).Value2 = 1620   # M Precio promedio ponderado
$ws.Cells.Item(41, 16).Value2 = 810    # P Precio $/Kg
